$wb = $excel.ActiveWorkbook

# Mapping of sheet name -> row -> new B-column ("ss") value, derived from the diff.
# Only column B (the standard-score lookup values) changes; column A (raw score) and
# the header row are untouched.
$changes = @{
    "K-Fall" = @{ 2=72; 3=73; 4=75; 5=77; 6=79; 7=80; 8=82; 9=84; 10=85; 11=87; 12=89; 15=93; 16=95; 17=97; 18=99; 19=100; 20=102; 21=104; 22=106; 23=108; 24=110; 25=112; 26=114; 27=117; 28=119; 29=123; 30=128; 31=130; 32=130; 33=130 }
    "K-Spring" = @{ 2=62; 3=64; 4=66; 5=68; 6=70; 7=71; 8=73; 9=75; 10=77; 11=79; 12=80; 13=82; 14=84; 15=85; 16=87; 17=89; 22=97; 23=99; 24=101; 25=103; 26=105; 27=107; 28=109; 29=111; 30=114; 31=116; 33=124; 34=130; 35=130 }
    "1-Fall" = @{ 2=49; 3=52; 4=54; 5=57; 6=59; 7=61; 8=64; 9=66; 10=68; 11=70; 12=72; 13=74; 14=75; 15=77; 16=79; 17=81; 18=83; 19=84; 20=86; 21=88; 24=93; 25=95; 26=97; 27=99; 28=101; 29=103; 30=106; 31=108; 32=111; 33=113; 35=123 }
    "1-Spring" = @{ 2=40; 3=40; 4=40; 5=42; 6=46; 7=49; 8=52; 9=54; 10=57; 11=59; 12=62; 13=64; 14=66; 15=68; 16=70; 17=72; 18=74; 19=76; 20=78; 21=80; 22=82; 23=84; 24=86; 25=88; 26=90; 27=92; 28=94; 29=96; 30=99; 31=101; 32=104; 33=107; 34=110; 35=115 }
    "2-Fall" = @{ 9=40; 10=42; 12=48; 13=51; 14=54; 15=57; 16=60; 17=62; 18=64; 19=67; 20=69; 21=71; 22=73; 23=76; 24=78; 25=80; 26=82; 29=89; 30=92; 31=94; 32=97; 33=100; 34=104; 35=110 }
    "2-Spring" = @{ 15=41; 17=48; 18=52; 19=55; 20=57; 21=60; 22=63; 23=65; 24=68; 25=70; 26=73; 27=75; 28=78; 29=81; 31=86; 32=89; 33=93; 34=97; 35=103 }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowValues = $changes[$sheetName]
    foreach ($row in $rowValues.Keys) {
        $ws.Cells.Item([int]$row, 2).Value = $rowValues[$row]
    }
}
